# edit.ps1
# 1) Update the cached "datetimeFigureOut" date placeholder text from
#    2/15/2024 -> 2/17/2024 on the slide master and every slide layout.
# 2) Split the title run on slide 1 so the leading word "Python" becomes
#    its own run ( " Python " ) and the remainder becomes
#    "Sudoku Program - Overview".

$p = $ppt.ActivePresentation

$oldDate = "2/15/2024"
$newDate = "2/17/2024"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDatePH = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePH = $true
            }
        } catch {
            $isDatePH = $false
        }
        if ($isDatePH -and $sh.HasTextFrame) {
            if ($sh.TextFrame.HasText) {
                $tr = $sh.TextFrame.TextRange
                if ($tr.Text -eq $oldDate) {
                    $tr.Text = $newDate
                }
            }
        }
    }
}

# Slide master
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout that hangs off the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# Title on slide 1: "Python Sudoku Program - Overview"
#   -> run1 " Python " + run2 "Sudoku Program - Overview"
$slide1 = $p.Slides.Item(1)
$titleShape = $slide1.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange

if ($titleRange.Text -eq "Python Sudoku Program - Overview") {
    $lead = $titleRange.Characters(1, 7)   # "Python "
    $lead.Text = " Python "
}
